# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / handoff / handback datetime
# stamps to reflect the newly generated report timestamps.

$wb = $excel.ActiveWorkbook

# Overview sheet: Latest HO Xliff Generate Date for d754c346-... row
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-09-07 01:02:08"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for d754c346-... row
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-09-07 01:01:58"
$wsZhCn.Range("K4").Value = "2016-09-07 01:02:33"

# de-de sheet: Correspond Handback DateTime for d754c346-... row
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K4").Value = "2016-09-07 01:02:41"
